$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in homework grades for row 13 (student 10)
$ws.Range("C13").Value = 5
$ws.Range("D13").Value = 5
$ws.Range("F13").Value = 5

# Update selection/active cell to G13
$ws.Range("G13").Select()
